# Rename header in "Weekly Quantity" sheet (B1): "Requested quantity" -> "Weekly_PO_Qty"
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# Rename header in "Monthly Trend" sheet (B1): "Requested quantity" -> "Monthly_PO_Qty"
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add a new worksheet "PO Forecast" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Data rows
$newSheet.Range("A2").Value = 45347.99999999999
$newSheet.Range("B2").Value = 243
$newSheet.Range("C2").Value = 162.0892017387136
$newSheet.Range("D2").Value = 317.659185463362
$newSheet.Range("A3").Value = 45361.99999999999
$newSheet.Range("B3").Value = 512
$newSheet.Range("C3").Value = 427.5827172033814
$newSheet.Range("D3").Value = 591.5067767975482
$newSheet.Range("A4").Value = 45396.99999999999
$newSheet.Range("B4").Value = 1185
$newSheet.Range("C4").Value = 1107.439657416279
$newSheet.Range("D4").Value = 1263.642582318307
$newSheet.Range("A5").Value = 45403.99999999999
$newSheet.Range("B5").Value = 1320
$newSheet.Range("C5").Value = 1240.507221927058
$newSheet.Range("D5").Value = 1398.35112386356
$newSheet.Range("A6").Value = 45410.99999999999
$newSheet.Range("B6").Value = 1455
$newSheet.Range("C6").Value = 1366.301813803967
$newSheet.Range("D6").Value = 1537.088912217351
$newSheet.Range("A7").Value = 45417.99999999999
$newSheet.Range("B7").Value = 1589
$newSheet.Range("C7").Value = 1511.71078079262
$newSheet.Range("D7").Value = 1675.662651465162
$newSheet.Range("A8").Value = 45424.99999999999
$newSheet.Range("B8").Value = 1724
$newSheet.Range("C8").Value = 1646.909625936938
$newSheet.Range("D8").Value = 1806.487100263024
$newSheet.Range("A9").Value = 45431.99999999999
$newSheet.Range("B9").Value = 1858
$newSheet.Range("C9").Value = 1779.072225640589
$newSheet.Range("D9").Value = 1944.317572718734
$newSheet.Range("A10").Value = 45438.99999999999
$newSheet.Range("B10").Value = 1993
$newSheet.Range("C10").Value = 1913.988173875882
$newSheet.Range("D10").Value = 2068.570528995293
$newSheet.Range("A11").Value = 45445.99999999999
$newSheet.Range("B11").Value = 2128
$newSheet.Range("C11").Value = 2040.935663273497
$newSheet.Range("D11").Value = 2206.048944456394
$newSheet.Range("A12").Value = 45452.99999999999
$newSheet.Range("B12").Value = 2262
$newSheet.Range("C12").Value = 2179.977367561123
$newSheet.Range("D12").Value = 2343.427768529511

# Apply header styling to match the other sheets (bold font, thin border,
# centered horizontal/vertical alignment)
$headerRange = $newSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Apply the shared date/time number format to the "ds" column (matches the
# formatting used for the date columns on the other sheets)
$dateRange = $newSheet.Range("A2:A12")
$dateRange.NumberFormat = "YYYY-MM-DD HH:MM:SS"
